# Validation with mothervessels performing tasks and good results?
#
# - task_compatibility: mother vessel (V4 / column E) is now marked as
#   compatible with every task (M1..M4), i.e. E2:E5 flip from 0 to 1.
# - Selection/active-sheet bookkeeping moves from "locations" to "vessels"
#   as the sheet left showing when the file was saved, and the selection
#   on task_compatibility moves to E5 (next to the values just edited).

$wb = $excel.ActiveWorkbook

# --- task_compatibility: allow the mother vessel (column E / V4) to
# perform all four tasks -------------------------------------------------
$wsTaskComp = $wb.Worksheets.Item("task_compatibility")
$wsTaskComp.Activate()
$wsTaskComp.Range("E2").Value = 1
$wsTaskComp.Range("E3").Value = 1
$wsTaskComp.Range("E4").Value = 1
$wsTaskComp.Range("E5").Value = 1
$wsTaskComp.Range("E5").Select()

# --- leave the workbook with "vessels" as the active / front-most sheet,
# with C5 selected, matching the saved view state of the edited file ----
$wsVessels = $wb.Worksheets.Item("vessels")
$wsVessels.Activate()
$wsVessels.Range("C5").Select()
